# Update "想去人数" (F) and "最低票价" (G) figures on both the 展览 and
# 全部类型 sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value  = 2937

    $ws.Range("F8").Value  = 1621

    $ws.Range("F9").Value  = 1604
    $ws.Range("G9").Value  = 55

    $ws.Range("F11").Value = 346

    $ws.Range("F15").Value = 23

    $ws.Range("F16").Value = 222

    $ws.Range("F23").Value = 344

    $ws.Range("F24").Value = 124

    $ws.Range("F27").Value = 1946

    $ws.Range("F31").Value = 151

    $ws.Range("F36").Value = 479
}
